$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 9 and 10 swap their identity data (index/prolificid/name/gender):
# Ankai (male, index 0) moves from row 9 to row 10.
# Annes (female, index 3) moves from row 10 to row 9.
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = "60bd88b8fc436774352f53b9"
$ws.Range("F9").Value = "Annes"
$ws.Range("G9").Value = "female"

$ws.Range("D10").Value = 0
$ws.Range("E10").Value = "5c27de12a2b00a00018b2c16"
$ws.Range("F10").Value = "Ankai"
$ws.Range("G10").Value = "male"

# Updated realeffort (H column) scores for all 24 data rows.
$ws.Range("H2").Value = 11.12356095231806
$ws.Range("H3").Value = 10.07657103797102
$ws.Range("H4").Value = 8.469824362969149
$ws.Range("H5").Value = 8.252431536799262
$ws.Range("H6").Value = 7.094657342882389
$ws.Range("H7").Value = 6.351992923050718
$ws.Range("H8").Value = 6.111943368614604
$ws.Range("H9").Value = 5.465857846036377
$ws.Range("H10").Value = 5.411049145544538
$ws.Range("H11").Value = 4.011565163053068
$ws.Range("H12").Value = 2.037811163075423
$ws.Range("H13").Value = 0.2234880371687213
$ws.Range("H14").Value = 8.487299773058886
$ws.Range("H15").Value = 8.17703893189862
$ws.Range("H16").Value = 7.389685286561035
$ws.Range("H17").Value = 7.23768709675885
$ws.Range("H18").Value = 6.343048961041616
$ws.Range("H19").Value = 6.040248810237871
$ws.Range("H20").Value = 5.132481009527025
$ws.Range("H21").Value = 3.39083832478149
$ws.Range("H22").Value = 3.318523148442117
$ws.Range("H23").Value = 2.176093329063497
$ws.Range("H24").Value = 1.474100674863718
$ws.Range("H25").Value = 0.05652031328557322
